$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update row 3: was Ravi / Off Duty -> becomes Meena / On Duty ---
$ws.Range("A3").Value = "Meena"
$ws.Range("B3").Value = 0.33333333333333331
$ws.Range("C3").Value = 0.66666666666666663
$ws.Range("D3").Value = "On Duty"

# --- Update row 4: was Meena / On Duty -> becomes Ravi / On Duty ---
$ws.Range("A4").Value = "Ravi"
$ws.Range("B4").Value = 0.70833333333333337
$ws.Range("C4").Value = 0.041666666666666664
$ws.Range("D4").Value = "On Duty"

# --- Add new row 5 for Sneha / Leave, copying formatting from row 4 ---
$ws.Range("A4:D4").Copy()
$ws.Range("A5:D5").PasteSpecial(-4122)

$ws.Range("A5").Value = "Sneha"
$ws.Range("B5").Value = 0.41666666666666669
$ws.Range("C5").Value = 0.75
$ws.Range("D5").Value = "Leave"

$ws.Range("A1:D5").Select() | Out-Null
